$d = $word.ActiveDocument

function Rename-Bookmark($oldName, $newName) {
    $bm = $d.Bookmarks.Item($oldName)
    $rng = $d.Range($bm.Start, $bm.Start)
    $bm.Delete()
    $d.Bookmarks.Add($newName, $rng)
}

# Figure 1
Rename-Bookmark "figure-1-figuresfig1-treetree-abbrev.pdf" "figure-1-figuresfig1-treefigure1.pdf"
$d.Content.Find.Execute(
    "Figure 1 (figures/fig1-tree/tree-abbrev.pdf)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Figure 1 (figures/fig1-tree/Figure1.pdf)", 2) | Out-Null

# Figure 2a(b) -> Figure 2
Rename-Bookmark "figure-2ab-figuresfig2-samplingcompleteness-2.svg-and-figuresfig2-samplingcompleteness-1.svg" "figure-2-figuresfig2-samplingfigure2.pdf"
$d.Content.Find.Execute(
    "Figure 2a(b) (figures/fig2-sampling/completeness-2.svg and figures/fig2-sampling/completeness-1.svg)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Figure 2 (figures/fig2-sampling/Figure2.pdf)", 2) | Out-Null

# Figure 4
Rename-Bookmark "figure-4-figuresfig4-seedsaci-aux-gh.pdf" "figure-4-figuresfig4-seedsfigure4.pdf"
$d.Content.Find.Execute(
    "Figure 4 (figures/fig4-seeds/acI-aux-GH.pdf)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Figure 4 (figures/fig4-seeds/Figure4.pdf)", 2) | Out-Null

# Figure 5
Rename-Bookmark "figure-5-figuresfig5-transportersaci-transporters.pdf" "figure-5-figuresfig5-transportersfigure5.pdf"
$d.Content.Find.Execute(
    "Figure 5 (figures/fig5-transporters/acI-transporters.pdf)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Figure 5 (figures/fig5-transporters/Figure5.pdf)", 2) | Out-Null

# Supplementary Figure 5
Rename-Bookmark "supplementary-figure-5-figuresfig1-treetree-full.pdf" "supplementary-figure-5-figuresfig1-treefigures5.pdf"
$d.Content.Find.Execute(
    "Supplementary Figure 5 (figures/fig1-tree/tree-full.pdf)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Supplementary Figure 5 (figures/fig1-tree/FigureS5.pdf)", 2) | Out-Null

Write-Output "done"
